# Apply the July 14, 2020 covid_disparities run results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - Maryland: previously-empty cells now populated with the
# successful run's values; status moves from an error message to "Success!".
$ws.Range("B17").Value = 44026
$ws.Range("B17").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C17").Value = 74260
$ws.Range("D17").Value = 3202
$ws.Range("E17").Value = 21525
$ws.Range("F17").Value = 1301
$ws.Range("G17").Value = 35.07
$ws.Range("H17").Value = 40.87

$ws.Range("K17").Value = 61384
$ws.Range("L17").Value = 3183

$ws.Range("O17").Value = "Success!"

# Row 36 - Iowa: total case count revised upward by one.
$ws.Range("C36").Value = 35866
